$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D (new "most recent year" column).
# This shifts the existing D:K data right to E:L, matching the target layout.
$ws.Columns("D").Insert()

# Copy number formats / styles from the (shifted) column E into the new column D
# so the new cells carry the same per-row styles (date format row 7/38/80, etc.)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the newest year's figures (FY ending 2018-12-31 = serial 43465).
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 881300
$ws.Range("D9").Value = 280900
$ws.Range("D10").Value = 600400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = -6100
$ws.Range("D15").Value = 309900
$ws.Range("D17").Value = 623200
$ws.Range("D18").Value = 258100
$ws.Range("D20").Value = 3900
$ws.Range("D21").Value = 571900
$ws.Range("D22").Value = 133400
$ws.Range("D23").Value = 128600
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 128600
$ws.Range("D27").Value = 115500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3900
$ws.Range("D33").Value = 115500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 115500

$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 146200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 129200
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 7800
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 98800
$ws.Range("D48").Value = 7787500
$ws.Range("D49").Value = 5200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 8261700
$ws.Range("D57").Value = 75100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 105700
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 4134000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 5859400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -935600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2402300
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 115500
$ws.Range("D83").Value = 309900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 433000
$ws.Range("D91").Value = -247500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -249600
$ws.Range("D96").Value = -169800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -213800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -30400
